# Edit script: applies the changes described by the commit diff to the CV document.
#
# Summary of changes:
#  1. Remove the "_GoBack" bookmark from its old location (after "... and ").
#  2. "showcased a machine learning proof-of-concept to" -> "showcased a LUIS machine learning proof-of-concept to"
#  3. "To complete AGNICO-EAGLES annual SOX certifications, audited ..." -> "Performed AGNICO-EAGLES annual SOX certifications: audited ..."
#  4. "12 requests of change in transaction monitoring" -> "12 change requests in transaction monitoring"
#  5. "... corporate actions workflow (more accurate dividends and splits)" ->
#     "... corporate actions workflow for accurate dividends and splits"
#     and the "_GoBack" bookmark is (re-)inserted right at the end of that bullet's text.

$d = $word.ActiveDocument

$wdFindContinue = 1
$wdReplaceOne = 1

function Replace-Text($find, $replace) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($find, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $replace, $wdReplaceOne)
    if (-not $ok) {
        Write-Output "WARNING: could not find text: $find"
    }
    return $ok
}

# 1. Remove the old "_GoBack" bookmark (was located right before "significant control deficiencies").
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2. "a machine learning proof-of-concept to" gains "LUIS "
Replace-Text " a machine learning proof-of-concept to " " a LUIS machine learning proof-of-concept to "

# 3. "To complete " -> "Performed "
Replace-Text "To complete " "Performed "

# 3b. "annual SOX certifications, a" -> "annual SOX certifications: a"
Replace-Text "annual SOX certifications, a" "annual SOX certifications: a"

# 4. "12 requests of change in " -> "12 change requests in "
Replace-Text "12 requests of change in " "12 change requests in "

# 5. "(more " -> "for "
Replace-Text "(more " "for "

# 5b. remove the trailing ")" after "accurate dividends and splits"
Replace-Text "accurate dividends and splits)" "accurate dividends and splits"

# 5c. Re-add the "_GoBack" bookmark right after "... splits" (end of that bullet's text,
#     immediately before the paragraph mark). Adding a bookmark exactly at "end of paragraph
#     text" position is unreliable, so as a workaround we temporarily insert a marker
#     character after the target point, anchor the bookmark just before that marker, then
#     remove the marker again -- the bookmark (a zero-length range) stays put.
$rng = $d.Content
$found = $rng.Find.Execute("accurate dividends and splits")
if ($found) {
    $endPos = $rng.End
    $marker = $d.Range($endPos, $endPos)
    $marker.InsertAfter("~")

    $bmPoint = $d.Range($endPos, $endPos)
    $d.Bookmarks.Add("_GoBack", $bmPoint)

    $markerRange = $d.Range($endPos, $endPos + 1)
    $markerRange.Delete()
} else {
    Write-Output "WARNING: could not find insertion point for _GoBack bookmark"
}

Write-Output "Done."
